# Apply the updated cryptos-list values scraped by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds values that look numeric (e.g. "60.896.67",
# "146.40", "0.999") but must stay plain text, exactly as authored in the
# sheet. Force the number format to Text before writing so Excel's COM
# layer doesn't silently coerce them into real numbers, then restore the
# cell's original (default/Normal) style so no stray style index lingers.
$priceUpdates = @{
    'D2' = '60.896.67'
    'D3' = '2.918.43'
    'D5' = '591.15'
    'D6' = '146.40'
    'D9' = '6.88'
    'D13' = '33.62'
    'D15' = '3.406.01'
    'D16' = '60.843.02'
    'D18' = '2.918.26'
    'D19' = '431.08'
    'D20' = '13.41'
    'D23' = '81.51'
    'D24' = '11.00'
    'D26' = '11.85'
    'D30' = '7.03'
    'D31' = '26.68'
    'D33' = '0.999'
    'D34' = '0.0₃0863'
    'D42' = '40.28'
    'D43' = '380.46'
    'D44' = '0.0343'
    'D45' = '2.694.21'
    'D46' = '133.33'
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = '@'
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = 'Normal'
}

# Coin names, links (row 44/45 were reordered: VeChain now before Maker)
# and the Volume(1h) percentages in column E are plain strings already,
# so a direct assignment is safe.
$plainUpdates = @{
    'E2' = '  +0.05%  '
    'E3' = '  -0.01%  '
    'E5' = '  +1.31%  '
    'E6' = '  +1.27%  '
    'E7' = '  -0.01%  '
    'E8' = '  +0.59%  '
    'E9' = '  +1.33%  '
    'E10' = '  -0.64%  '
    'E11' = '  -1.49%  '
    'E12' = '  -0.09%  '
    'E13' = '  -0.05%  '
    'E14' = '  +0.11%  '
    'E15' = '  +0.04%  '
    'E16' = '  +0.00%  '
    'E17' = '  -0.73%  '
    'E18' = '  -0.07%  '
    'E19' = '  +0.10%  '
    'E20' = '  -1.87%  '
    'E21' = '  -0.83%  '
    'E22' = '  -1.06%  '
    'E23' = '  +1.13%  '
    'E24' = '  +1.28%  '
    'E25' = '  -0.76%  '
    'E26' = '  -0.49%  '
    'E27' = '  -0.03%  '
    'E28' = '  +5.29%  '
    'E29' = '  -0.11%  '
    'E30' = '  -2.88%  '
    'E31' = '  +0.36%  '
    'E32' = '  +1.85%  '
    'E33' = '  -0.12%  '
    'E34' = '  -1.23%  '
    'E35' = '  -0.04%  '
    'E36' = '  -0.45%  '
    'E37' = '  +0.49%  '
    'E38' = '  -1.33%  '
    'E39' = '  -4.37%  '
    'E40' = '  -1.57%  '
    'E41' = '  -3.14%  '
    'E42' = '  -2.60%  '
    'E43' = '  +0.65%  '
    'B44' = 'VeChain'
    'C44' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'E44' = '  -1.62%  '
    'B45' = 'Maker'
    'C45' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'E45' = '  +0.50%  '
    'E46' = '  +0.51%  '
    'E48' = '  -2.48%  '
    'E49' = '  -0.58%  '
    'E50' = '  -2.85%  '
    'E51' = '  -0.03%  '
}

foreach ($addr in $plainUpdates.Keys) {
    $ws.Range($addr).Value = $plainUpdates[$addr]
}

Write-Output "Applied $($priceUpdates.Count + $plainUpdates.Count) cell updates"
